{"js": "// 1) \"shipper\" -> \"seller\" inside the {{shipper_address}} merge tag.\nconst shipperResults = context.document.body.search(\"shipper\", { matchCase: true });\nshipperResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < shipperResults.items.length; i++) {\n  shipperResults.items[i].insertText(\"seller\", Word.InsertLocation.replace);\n}\n\n// 2) \"consignee\" -> \"buyer\" inside the {{consignee_address}} merge tag.\nconst consigneeResults = context.document.body.search(\"consignee\", { matchCase: true });\nconsigneeResults.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < consigneeResults.items.length; i++) {\n  consigneeResults.items[i].insertText(\"buyer\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n\n// 3) Remove the two \"\u0418\u0434\u0435\u043d\u0442\u0438\u0444\u0438\u043a\u0430\u0442\u043e\u0440 \u0433\u043e\u0441\u0443\u0434\u0430\u0440\u0441\u0442\u0432\u0435\u043d\u043d\u043e\u0433\u043e \u043a\u043e\u043d\u0442\u0440\u0430\u043a\u0442\u0430 ... (8)\" paragraphs.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst toDelete = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"\u0418\u0434\u0435\u043d\u0442\u0438\u0444\u0438\u043a\u0430\u0442\u043e\u0440\") === 0 || text.indexOf(\"\u043d\u0430\u043b\u0438\u0447\u0438\u0438)\") === 0) {\n    toDelete.push(paragraphs.items[i]);\n  }\n}\nfor (let i = 0; i < toDelete.length; i++) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n\n// 4) Shrink the page margins to 0.5in (720 twips = 36pt) on every side.\nconst section = context.document.sections.getFirst();\nconst pageSetup = section.pageSetup;\npageSetup.topMargin = 36;\npageSetup.bottomMargin = 36;\npageSetup.leftMargin = 36;\npageSetup.rightMargin = 36;\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"shipper\" -> \"seller\" inside the {{shipper_address}} merge tag.\n$find1 = $d.Content\n$find1.Find.ClearFormatting()\n$find1.Find.Replacement.ClearFormatting()\n$find1.Find.Text = \"shipper\"\n$find1.Find.Replacement.Text = \"seller\"\n$find1.Find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 2) \"consignee\" -> \"buyer\" inside the {{consignee_address}} merge tag.\n$find2 = $d.Content\n$find2.Find.ClearFormatting()\n$find2.Find.Replacement.ClearFormatting()\n$find2.Find.Text = \"consignee\"\n$find2.Find.Replacement.Text = \"buyer\"\n$find2.Find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# 3) Remove the two \"\u0418\u0434\u0435\u043d\u0442\u0438\u0444\u0438\u043a\u0430\u0442\u043e\u0440 \u0433\u043e\u0441\u0443\u0434\u0430\u0440\u0441\u0442\u0432\u0435\u043d\u043d\u043e\u0433\u043e \u043a\u043e\u043d\u0442\u0440\u0430\u043a\u0442\u0430 ... (8)\" paragraphs.\n$toDelete = @()\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text\n  if ($t.StartsWith(\"\u0418\u0434\u0435\u043d\u0442\u0438\u0444\u0438\u043a\u0430\u0442\u043e\u0440\") -or $t.StartsWith(\"\u043d\u0430\u043b\u0438\u0447\u0438\u0438)\")) {\n    $toDelete += $i\n  }\n}\nfor ($j = $toDelete.Count - 1; $j -ge 0; $j--) {\n  $d.Paragraphs.Item($toDelete[$j]).Range.Delete() | Out-Null\n}\n\n# 4) Shrink the page margins to 0.5in (720 twips = 36pt) on every side.\n$ps = $d.PageSetup\n$ps.TopMargin = 36\n$ps.BottomMargin = 36\n$ps.LeftMargin = 36\n$ps.RightMargin = 36\n"}
